$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$ws.Range("H17").Value = 828996.8
$ws.Range("J17").Value = 879659.9
$ws.Range("L17").Value = 2638979.7
$ws.Range("N17").Value = -2639315.7

# Row 62 (ALC)
$ws.Range("H62").Value = 3183023.5
$ws.Range("I62").Value = 5563073.5
$ws.Range("J62").Value = 9623.866
$ws.Range("K62").Value = 5563073.5
$ws.Range("L62").Value = 9623.866
$ws.Range("M62").Value = -5562449.5
$ws.Range("N62").Value = -10871.866

# Row 65 (ALC)
$ws.Range("H65").Value = 3183023.5
$ws.Range("I65").Value = 5563073.5
$ws.Range("J65").Value = 9623.866
$ws.Range("K65").Value = 27815367.5
$ws.Range("L65").Value = 48119.33
$ws.Range("M65").Value = -27812247.5
$ws.Range("N65").Value = -54359.33

# Row 121 (ALC)
$ws.Range("H121").Value = 851.25
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 1001.6667
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 3005.0001
$ws.Range("M121").Value = 547
$ws.Range("N121").Value = -6499.0001

# Row 132 (ALC)
$ws.Range("H132").Value = 195396.45
$ws.Range("I132").Value = 248897.38
$ws.Range("J132").Value = 8143.2856
$ws.Range("K132").Value = 746692.14
$ws.Range("L132").Value = 24429.8568
$ws.Range("M132").Value = -744162.14
$ws.Range("N132").Value = -29489.8568

# Row 137 (ALC)
$ws.Range("H137").Value = 19231758
$ws.Range("I137").Value = 26316418
$ws.Range("J137").Value = 1965.5
$ws.Range("K137").Value = 78949254
$ws.Range("L137").Value = 5896.5
$ws.Range("M137").Value = -78946704
$ws.Range("N137").Value = -10996.5

# Row 138 (ALC)
$ws.Range("H138").Value = 3906996
$ws.Range("I138").Value = 859217.75
$ws.Range("J138").Value = 7578183.5
$ws.Range("K138").Value = 2577653.25
$ws.Range("L138").Value = 22734550.5
$ws.Range("M138").Value = -2572513.25
$ws.Range("N138").Value = -22744830.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 15081.917
$ws.Range("I32").Value = 2867.87
$ws.Range("J32").Value = 149436.42
$ws.Range("K32").Value = 2867.87
$ws.Range("L32").Value = 149436.42
$ws.Range("M32").Value = -2580.87
$ws.Range("N32").Value = -150010.42

# Row 45 (ARM)
$ws.Range("H45").Value = 859.2857
$ws.Range("I45").Value = 809.1667
$ws.Range("K45").Value = 809.1667
$ws.Range("M45").Value = -432.1667

# Row 61 (ARM)
$ws.Range("H61").Value = 2489.1555
$ws.Range("I61").Value = 1811.027
$ws.Range("J61").Value = 5625.5
$ws.Range("K61").Value = 1811.027
$ws.Range("L61").Value = 5625.5
$ws.Range("M61").Value = -1599.027
$ws.Range("N61").Value = -6049.5

# Row 74 (ARM)
$ws.Range("H74").Value = 4148.814
$ws.Range("I74").Value = 1111.5938
$ws.Range("K74").Value = 1111.5938
$ws.Range("M74").Value = -237.5938000000001

# Row 77 (ARM)
$ws.Range("H77").Value = 4148.814
$ws.Range("I77").Value = 1111.5938
$ws.Range("K77").Value = 5557.969000000001
$ws.Range("M77").Value = -1189.969000000001

# Row 102 (ARM)
$ws.Range("H102").Value = 6222.2
$ws.Range("I102").Value = 6527.75
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 6527.75
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -4905.75
$ws.Range("N102").Value = -8244

# Row 110 (ARM)
$ws.Range("H110").Value = 1710.7333
$ws.Range("I110").Value = 852.2
$ws.Range("J110").Value = 2140
$ws.Range("K110").Value = 852.2
$ws.Range("L110").Value = 2140
$ws.Range("M110").Value = 1192.8
$ws.Range("N110").Value = -6230

# Row 132 (ARM)
$ws.Range("H132").Value = 2281.5576
$ws.Range("I132").Value = 2010.0698
$ws.Range("J132").Value = 3578.6667
$ws.Range("K132").Value = 6030.2094
$ws.Range("L132").Value = 10736.0001
$ws.Range("M132").Value = -3500.2094
$ws.Range("N132").Value = -15796.0001

# Row 136 (ARM)
$ws.Range("H136").Value = 2489.1555
$ws.Range("I136").Value = 1811.027
$ws.Range("J136").Value = 5625.5
$ws.Range("K136").Value = 5433.081
$ws.Range("L136").Value = 16876.5
$ws.Range("M136").Value = -2883.081
$ws.Range("N136").Value = -21976.5

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (BSM)
$ws.Range("H134").Value = 14494348
$ws.Range("I134").Value = 15874219
$ws.Range("K134").Value = 47622657
$ws.Range("M134").Value = -47620122

$ws = $wb.Worksheets.Item("CRP")
# Row 9 (CRP)
$ws.Range("H9").Value = 88888
$ws.Range("J9").Value = 88888
$ws.Range("L9").Value = 88888
$ws.Range("N9").Value = -89224

# Row 16 (CRP)
$ws.Range("H16").Value = 24760.953
$ws.Range("I16").Value = 42617
$ws.Range("J16").Value = 952.8889
$ws.Range("K16").Value = 42617
$ws.Range("L16").Value = 952.8889
$ws.Range("M16").Value = -42330
$ws.Range("N16").Value = -1526.8889

# Row 31 (CRP)
$ws.Range("H31").Value = 1772.4036
$ws.Range("I31").Value = 1084.7949
$ws.Range("J31").Value = 3262.2222
$ws.Range("K31").Value = 1084.7949
$ws.Range("L31").Value = 3262.2222
$ws.Range("M31").Value = -789.7949000000001
$ws.Range("N31").Value = -3852.2222

# Row 34 (CRP)
$ws.Range("H34").Value = 1772.4036
$ws.Range("I34").Value = 1084.7949
$ws.Range("J34").Value = 3262.2222
$ws.Range("K34").Value = 1084.7949
$ws.Range("L34").Value = 3262.2222
$ws.Range("M34").Value = -882.7949000000001
$ws.Range("N34").Value = -3666.2222

# Row 58 (CRP)
$ws.Range("H58").Value = 2006.1724
$ws.Range("I58").Value = 1272.7727
$ws.Range("J58").Value = 4311.143
$ws.Range("K58").Value = 1272.7727
$ws.Range("L58").Value = 4311.143
$ws.Range("M58").Value = -1069.7727
$ws.Range("N58").Value = -4717.143

# Row 99 (CRP)
$ws.Range("H99").Value = 2981722.8
$ws.Range("I99").Value = 4815136.5
$ws.Range("K99").Value = 4815136.5
$ws.Range("M99").Value = -4813638.5

# Row 113 (CRP)
$ws.Range("H113").Value = 24760.953
$ws.Range("I113").Value = 42617
$ws.Range("J113").Value = 952.8889
$ws.Range("K113").Value = 42617
$ws.Range("L113").Value = 952.8889
$ws.Range("M113").Value = -40447
$ws.Range("N113").Value = -5292.8889

# Row 126 (CRP)
$ws.Range("H126").Value = 2981722.8
$ws.Range("I126").Value = 4815136.5
$ws.Range("K126").Value = 14445409.5
$ws.Range("M126").Value = -14442939.5

# Row 132 (CRP)
$ws.Range("H132").Value = 1752.0968
$ws.Range("I132").Value = 1424.5714
$ws.Range("J132").Value = 4809
$ws.Range("K132").Value = 4273.7142
$ws.Range("L132").Value = 14427
$ws.Range("M132").Value = -1743.7142
$ws.Range("N132").Value = -19487.0005

# Row 136 (CRP)
$ws.Range("H136").Value = 2006.1724
$ws.Range("I136").Value = 1272.7727
$ws.Range("J136").Value = 4311.143
$ws.Range("K136").Value = 3818.3181
$ws.Range("L136").Value = 12933.429
$ws.Range("M136").Value = -1268.3181
$ws.Range("N136").Value = -18033.429

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (CUL)
$ws.Range("H131").Value = 6291005
$ws.Range("I131").Value = 504.14285
$ws.Range("J131").Value = 7248255.5
$ws.Range("K131").Value = 1512.42855
$ws.Range("L131").Value = 21744766.5
$ws.Range("M131").Value = 3527.57145
$ws.Range("N131").Value = -21754846.5

$ws = $wb.Worksheets.Item("LTW")
# Row 132 (LTW)
$ws.Range("H132").Value = 2023.2206
$ws.Range("I132").Value = 1271.9584
$ws.Range("J132").Value = 3826.25
$ws.Range("K132").Value = 3815.8752
$ws.Range("L132").Value = 11478.75
$ws.Range("M132").Value = -1285.8752
$ws.Range("N132").Value = -16538.75

# Row 136 (LTW)
$ws.Range("H136").Value = 5131.3
$ws.Range("I136").Value = 3145.889
$ws.Range("J136").Value = 23000
$ws.Range("K136").Value = 9437.667000000001
$ws.Range("L136").Value = 69000
$ws.Range("M136").Value = -6887.667000000001
$ws.Range("N136").Value = -74100

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR)
$ws.Range("H132").Value = 2337.7876
$ws.Range("I132").Value = 2308.2307
$ws.Range("J132").Value = 2465.8667
$ws.Range("K132").Value = 6924.6921
$ws.Range("L132").Value = 7397.6001
$ws.Range("M132").Value = -4394.6921
$ws.Range("N132").Value = -12457.6001

# Row 136 (WVR)
$ws.Range("H136").Value = 28481.756
$ws.Range("I136").Value = 42443.125
$ws.Range("J136").Value = 2706.923
$ws.Range("K136").Value = 127329.375
$ws.Range("L136").Value = 8120.768999999999
$ws.Range("M136").Value = -124779.375
$ws.Range("N136").Value = -13220.769
